# Updated cryptos list on Mon Feb 27 01:13:23 UTC 2023 with GitHub Actions
#
# Applies the new price / volume(1h) values to the cryptos worksheet, and
# accounts for the rows that swapped order (Filecoin <-> WrappedliquidstakedEther2.0,
# InternetComputer(DFINITY) <-> Algorand, Aptos <-> EnergySwap, PancakeSwap <-> Frax).
#
# Values are written with a leading apostrophe so Excel stores them as literal
# text (matching the original inlineStr / text cells) instead of auto-converting
# strings such as "1.001" or "23.478.62" into numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'23.478.62"
$ws.Range("E2").Value = "'  +1.55%  "
$ws.Range("D3").Value = "'1.636.36"
$ws.Range("E3").Value = "'  +2.76%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "'  +0.21%  "
$ws.Range("D5").Value = "'308.15"
$ws.Range("E5").Value = "'  +2.20%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "'  +0.09%  "
$ws.Range("D7").Value = "'0.3775"
$ws.Range("E7").Value = "'  -0.04%  "
$ws.Range("D8").Value = "'52.94"
$ws.Range("E8").Value = "'  +4.15%  "
$ws.Range("D9").Value = "'0.3665"
$ws.Range("E9").Value = "'  +2.01%  "
$ws.Range("D10").Value = "'1.269"
$ws.Range("E10").Value = "'  +2.26%  "
$ws.Range("D11").Value = "'0.08184"
$ws.Range("E11").Value = "'  +1.37%  "
$ws.Range("D12").Value = "'1.002"
$ws.Range("E12").Value = "'  +0.20%  "
$ws.Range("D13").Value = "'23.05"
$ws.Range("E13").Value = "'  +3.27%  "
$ws.Range("D14").Value = "'6.661"
$ws.Range("E14").Value = "'  +2.12%  "
$ws.Range("D15").Value = "'0.00001279"
$ws.Range("E15").Value = "'  +3.03%  "
$ws.Range("D16").Value = "'7.451"
$ws.Range("E16").Value = "'  +1.53%  "
$ws.Range("D17").Value = "'1.638.61"
$ws.Range("E17").Value = "'  +2.97%  "
$ws.Range("D18").Value = "'94.84"
$ws.Range("E18").Value = "'  +1.73%  "
$ws.Range("D19").Value = "'0.06955"
$ws.Range("E19").Value = "'  +2.19%  "
$ws.Range("D20").Value = "'18.32"
$ws.Range("E20").Value = "'  +2.36%  "
$ws.Range("D21").Value = "'6.584"
$ws.Range("E21").Value = "'  +1.84%  "
$ws.Range("D22").Value = "'0.9997"
$ws.Range("E22").Value = "'  -0.09%  "
$ws.Range("D23").Value = "'23.483.62"
$ws.Range("E23").Value = "'  +1.55%  "
$ws.Range("D24").Value = "'12.84"
$ws.Range("E24").Value = "'  -0.04%  "
$ws.Range("D25").Value = "'3.085"
$ws.Range("E25").Value = "'  +4.64%  "
$ws.Range("D26").Value = "'2.420"
$ws.Range("E26").Value = "'  +1.56%  "
$ws.Range("D27").Value = "'21.31"
$ws.Range("E27").Value = "'  +1.91%  "
$ws.Range("D28").Value = "'150.89"
$ws.Range("E28").Value = "'  +1.18%  "
$ws.Range("D29").Value = "'5.326"
$ws.Range("E29").Value = "'  +1.63%  "
$ws.Range("D30").Value = "'136.07"
$ws.Range("E30").Value = "'  +2.48%  "
$ws.Range("D31").Value = "'2.389"
$ws.Range("E31").Value = "'  +0.37%  "
$ws.Range("B32").Value = "'WrappedliquidstakedEther2.0"
$ws.Range("C32").Value = "'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D32").Value = "'1.819.42"
$ws.Range("E32").Value = "'  +3.00%  "
$ws.Range("B33").Value = "'Filecoin"
$ws.Range("C33").Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "'6.815"
$ws.Range("E33").Value = "'  +1.93%  "
$ws.Range("D34").Value = "'0.9746"
$ws.Range("E34").Value = "'  +0.62%  "
$ws.Range("D35").Value = "'0.02836"
$ws.Range("E35").Value = "'  +5.95%  "
$ws.Range("D36").Value = "'10.43"
$ws.Range("E36").Value = "'  +2.81%  "
$ws.Range("D37").Value = "'0.07417"
$ws.Range("E37").Value = "'  -1.08%  "
$ws.Range("B38").Value = "'Algorand"
$ws.Range("C38").Value = "'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D38").Value = "'0.2544"
$ws.Range("E38").Value = "'  +1.93%  "
$ws.Range("B39").Value = "'InternetComputer(DFINITY)"
$ws.Range("C39").Value = "'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D39").Value = "'6.187"
$ws.Range("E39").Value = "'  +1.54%  "
$ws.Range("D40").Value = "'0.08902"
$ws.Range("E40").Value = "'  +1.45%  "
$ws.Range("D41").Value = "'1.384"
$ws.Range("E41").Value = "'  +1.91%  "
$ws.Range("D42").Value = "'0.7117"
$ws.Range("E42").Value = "'  +0.83%  "
$ws.Range("B43").Value = "'EnergySwap"
$ws.Range("C43").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").Value = "'16.32"
$ws.Range("E43").Value = "'  +9.41%  "
$ws.Range("B44").Value = "'Aptos"
$ws.Range("C44").Value = "'https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D44").Value = "'12.53"
$ws.Range("E44").Value = "'  +2.63%  "
$ws.Range("D45").Value = "'0.6560"
$ws.Range("E45").Value = "'  +1.52%  "
$ws.Range("D46").Value = "'2.353"
$ws.Range("E46").Value = "'  +3.31%  "
$ws.Range("B47").Value = "'Frax"
$ws.Range("C47").Value = "'https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D47").Value = "'0.9996"
$ws.Range("E47").Value = "'  +0.05%  "
$ws.Range("B48").Value = "'PancakeSwap"
$ws.Range("C48").Value = "'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D48").Value = "'4.040"
$ws.Range("E48").Value = "'  +1.29%  "
$ws.Range("D49").Value = "'0.08037"
$ws.Range("E49").Value = "'  +1.93%  "
$ws.Range("D50").Value = "'129.92"
$ws.Range("E50").Value = "'  -1.08%  "
$ws.Range("D51").Value = "'1.213"
$ws.Range("E51").Value = "'  +0.51%  "
